$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "mtest_status"
$ws.Range("B7").Value = "7C0"
$ws.Range("C7").Value = "M"
$ws.Range("D7").Value = "R/W"
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = "standard"
$ws.Range("G7").Value = "Custom Register for handelling test success/failure in simulation"

$ws.Range("A9").Select()
